$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values must be written with a leading apostrophe so Excel
# keeps storing them as text (matching the original inline-string cells)
# rather than re-typing them as numbers; the style is then reset to the
# original "Normal" so no stray number-format/quote-prefix style lingers.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "245.16"   # D2: 245.35 -> 245.16
Set-TextValue 2 7 "22"   # G2: 21 -> 22
Set-TextValue 3 4 "24.20"   # D3: 24.18 -> 24.20
Set-TextValue 3 7 "22"   # G3: 21 -> 22
Set-TextValue 4 4 "5.281"   # D4: 5.278 -> 5.281
Set-TextValue 4 7 "22"   # G4: 21 -> 22
Set-TextValue 5 7 "22"   # G5: 21 -> 22
Set-TextValue 6 4 "6.455"   # D6: 6.458 -> 6.455
Set-TextValue 6 7 "22"   # G6: 21 -> 22
Set-TextValue 7 4 "3.133"   # D7: 3.143 -> 3.133
Set-TextValue 7 7 "22"   # G7: 21 -> 22
Set-TextValue 8 4 "0.8153"   # D8: 0.8168 -> 0.8153
Set-TextValue 8 7 "22"   # G8: 21 -> 22
Set-TextValue 9 4 "0.8502"   # D9: 0.8522 -> 0.8502
Set-TextValue 9 7 "22"   # G9: 21 -> 22
Set-TextValue 10 4 "0.1355"   # D10: 0.1359 -> 0.1355
Set-TextValue 10 7 "22"   # G10: 21 -> 22
Set-TextValue 11 4 "0.06943"   # D11: 0.06940 -> 0.06943
Set-TextValue 11 7 "22"   # G11: 21 -> 22
Set-TextValue 12 4 "0.03139"   # D12: 0.03126 -> 0.03139
Set-TextValue 12 7 "22"   # G12: 21 -> 22
Set-TextValue 13 4 "0.02907"   # D13: 0.02897 -> 0.02907
Set-TextValue 13 7 "22"   # G13: 21 -> 22
Set-TextValue 14 4 "0.09383"   # D14: 0.09380 -> 0.09383
Set-TextValue 14 7 "22"   # G14: 21 -> 22
Set-TextValue 15 4 "3.756"   # D15: 3.745 -> 3.756
Set-TextValue 15 7 "22"   # G15: 21 -> 22
Set-TextValue 16 7 "22"   # G16: 21 -> 22
Set-TextValue 17 4 "0.04668"   # D17: 0.04690 -> 0.04668
Set-TextValue 17 7 "22"   # G17: 21 -> 22
Set-TextValue 18 4 "0.0006011"   # D18: 0.0005989 -> 0.0006011
Set-TextValue 18 7 "22"   # G18: 21 -> 22
Set-TextValue 19 4 "0.006250"   # D19: 0.006217 -> 0.006250
Set-TextValue 19 7 "22"   # G19: 21 -> 22
Set-TextValue 20 4 "0.001237"   # D20: 0.001235 -> 0.001237
Set-TextValue 20 7 "22"   # G20: 21 -> 22
Set-TextValue 21 4 "0.004613"   # D21: 0.004616 -> 0.004613
Set-TextValue 21 7 "22"   # G21: 21 -> 22
Set-TextValue 22 4 "0.00006900"   # D22: 0.00006899 -> 0.00006900
Set-TextValue 22 7 "22"   # G22: 21 -> 22
Set-TextValue 23 4 "3.498"   # D23: 3.500 -> 3.498
Set-TextValue 23 7 "22"   # G23: 21 -> 22
Set-TextValue 24 4 "2.150"   # D24: 2.148 -> 2.150
Set-TextValue 24 7 "22"   # G24: 21 -> 22
Set-TextValue 25 4 "0.3192"   # D25: 0.3193 -> 0.3192
Set-TextValue 25 7 "22"   # G25: 21 -> 22
Set-TextValue 26 7 "22"   # G26: 21 -> 22
Set-TextValue 27 7 "22"   # G27: 21 -> 22
Set-TextValue 28 4 "0.0002332"   # D28: 0.0002331 -> 0.0002332
Set-TextValue 28 7 "22"   # G28: 21 -> 22
Set-TextValue 29 7 "22"   # G29: 21 -> 22
Set-TextValue 30 7 "22"   # G30: 21 -> 22
Set-TextValue 31 7 "22"   # G31: 21 -> 22
Set-TextValue 32 7 "22"   # G32: 21 -> 22
Set-TextValue 33 7 "22"   # G33: 21 -> 22
Set-TextValue 34 7 "22"   # G34: 21 -> 22
Set-TextValue 35 7 "22"   # G35: 21 -> 22
Set-TextValue 36 7 "22"   # G36: 21 -> 22
Set-TextValue 37 7 "22"   # G37: 21 -> 22
Set-TextValue 38 7 "22"   # G38: 21 -> 22
Set-TextValue 39 7 "22"   # G39: 21 -> 22
Set-TextValue 40 4 "0.03634"   # D40: 0.03644 -> 0.03634
Set-TextValue 40 7 "22"   # G40: 21 -> 22
Set-TextValue 41 4 "0.006240"   # D41: 0.006249 -> 0.006240
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"   # E41: 40KickTokenKICKBestin24h -> 40KickTokenKICK
Set-TextValue 41 7 "22"   # G41: 21 -> 22
Set-TextValue 42 4 "0.1052"   # D42: 0.1053 -> 0.1052
Set-TextValue 42 7 "22"   # G42: 21 -> 22
Set-TextValue 43 4 "0.002760"   # D43: 0.003400 -> 0.002760
Set-TextValue 43 7 "22"   # G43: 21 -> 22
Set-TextValue 44 4 "0.008407"   # D44: 0.008422 -> 0.008407
$ws.Cells.Item(44, 5).Value = "43LocalTradersLCTBestin24h"   # E44: 43LocalTradersLCT -> 43LocalTradersLCTBestin24h
Set-TextValue 44 7 "22"   # G44: 21 -> 22
Set-TextValue 45 7 "22"   # G45: 21 -> 22
Set-TextValue 46 7 "22"   # G46: 21 -> 22
Set-TextValue 47 4 "0.3701"   # D47: 0.3699 -> 0.3701
Set-TextValue 47 7 "22"   # G47: 21 -> 22
Set-TextValue 48 4 "0.002282"   # D48: 0.002283 -> 0.002282
Set-TextValue 48 7 "22"   # G48: 21 -> 22
Set-TextValue 49 7 "22"   # G49: 21 -> 22
Set-TextValue 50 7 "22"   # G50: 21 -> 22
Set-TextValue 51 7 "22"   # G51: 21 -> 22
